$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original data block (rows 2-7) gets duplicated twice more, filling
# rows 8-13 and then 14-19 with the same product rows, pushing nothing else
# around - the trailing formatted-but-empty row stays at row 23 (its
# `spans` just widens because the sheet's used range now reaches column A-D
# on every row).
$data = @(
    @("Alicate Universal 8`"", 35,   "https://i.imgur.com/zanG3Hx.jpg",  "produtos"),
    @("Torquês Armador",       45,   "https://i.imgur.com/HpwMEIn.jpeg", "produtos"),
    @("Estilete 1",            11.5, "https://i.imgur.com/ZZtUnso.jpeg", "produtos"),
    @("Estilete 2",            9.5,  "https://i.imgur.com/XePLgzW.jpeg", "produtos"),
    @("Alicate Universal 8`"", 30,   "https://i.imgur.com/zanG3Hx.jpg",  "ofertas"),
    @("Torquês Armador",       40,   "https://i.imgur.com/HpwMEIn.jpeg", "ofertas")
)

$destRow = 8
for ($copy = 0; $copy -lt 2; $copy++) {
    $srcRow = 2
    foreach ($row in $data) {
        # Carry over the number-format (column B) and wrap/valign (column D)
        # styling from the matching source row before writing values, so the
        # new cells share the same style indices instead of minting new ones.
        $ws.Cells.Item($srcRow, 2).Copy()
        $ws.Cells.Item($destRow, 2).PasteSpecial(-4122)
        $ws.Cells.Item($srcRow, 4).Copy()
        $ws.Cells.Item($destRow, 4).PasteSpecial(-4122)

        $ws.Cells.Item($destRow, 1).Value = $row[0]
        $ws.Cells.Item($destRow, 2).Value = $row[1]
        $ws.Cells.Item($destRow, 3).Value = $row[2]
        $ws.Cells.Item($destRow, 4).Value = $row[3]

        $srcRow++
        $destRow++
    }
}

$excel.CutCopyMode = $false

$ws.Range("A24").Select()
